$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") switches its table
#    style from {F697415C-5CD2-45FC-90AB-E7705E0BA30D} to
#    {A2B2CD30-E7DD-48F8-B155-982A529CF5BC}.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A2B2CD30-E7DD-48F8-B155-982A529CF5BC}")
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme colour palette is swapped from the "Integral" / "Red
#    Violet" palette to the stock "Office Theme" / "Office" palette (the
#    font scheme and format scheme are identical between the two themes --
#    only the twelve theme colours, in clrScheme order, actually change).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$colors = $slide1.ThemeColorScheme

# clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeTheme[$i - 1]
}
